$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply number formats (creates the new custom numFmt "0.000" and the
#     built-in "0.00" cellXfs entries, in the same order as the target file:
#     the custom "0.000" format is referenced first so it gets the lower
#     cellXfs index). ---

# "Geom opt (Hartree)" column E gets a 3-decimal format on a handful of rows.
$ws.Range("E25").NumberFormat = "0.000"
$ws.Range("E27:E31").NumberFormat = "0.000"

# "Geom opt (kcal/mol)" and "Single pt (kcal/mol)" columns (F and H) get a
# 2-decimal format for rows 6-31.
$ws.Range("F6:F31").NumberFormat = "0.00"
$ws.Range("H6:H31").NumberFormat = "0.00"

# --- Update the recalculated energy values in the last three data rows ---

# Row 29 (Active / fGln123 / Val51 / INT1)
$ws.Range("E29").Value = -5709.9439400000001
$ws.Range("F29").Value = 10.272338700000001
$ws.Range("G29").Value = -5713.9138999999996
$ws.Range("H29").Value = 10.960717199999999

# Row 30 (Active / fGln123 / Val51 / TS1)
$ws.Range("E30").Value = -5709.9393369999998
$ws.Range("F30").Value = 13.160767229999999
$ws.Range("G30").Value = -5713.9032999999999
$ws.Range("H30").Value = 17.579065100000001

# Row 31 (Active / fGln123 / Val51 / INT2)
$ws.Range("E31").Value = -5709.9216059999999
$ws.Range("F31").Value = 24.286999999999999
$ws.Range("G31").Value = -5713.9026000000003
$ws.Range("H31").Value = 18.021459700000001

# --- Update the active cell / selection shown when the sheet is opened ---
$ws.Range("K16").Select()
